# Add pickle load functionality
# Fills in the Module Name / Line Numbers columns for the rows that were
# previously left blank under the "Child class-2.1" (pickle) and
# "Vector operators" sections, and updates the sheet's view/selection
# state to reflect where the author was last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 31: "Read data from a pickle file"
$ws.Range("C31").Value = "Departments.py"
$ws.Range("D31").Value = 484

# Row 32: "Utilize configuration constants"
$ws.Range("C32").Value = "Departments.py"
$ws.Range("D32").Value = 219

# Row 33: "Visualize"
$ws.Range("C33").Value = "Departments.py"
$ws.Range("D33").Value = 313

# Row 41: "Vector operators: Display, export"
$ws.Range("C41").Value = "Departments.py"
$ws.Range("D41").Value = 378

# Scroll / selection state: the author ended up with the view scrolled to
# show row 12 onward, with the last-edited cell (D41) selected.
$excel.ActiveWindow.ScrollRow = 12
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D41").Select()
